# Update column C ("Förändrad") for rows 2-43 from 45765 (2025-04-18)
# to 45766 (2025-04-19) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45765) {
        $cell.Value2 = 45766
    }
}
